$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = "c76eb7bc1a2e6e67034fc7d750762de9"
$ws.Range("B11").Value = "4655c8946a46d00e5c34817b6b9e389d"
$ws.Range("B17").Value = "6d8ffd1d66c53c710be8772851e1d28a"
$ws.Range("B24").Value = "f922ed9e447644263a1a97de707e5cb8"
$ws.Range("B29").Value = "b260b6fab631e6ee7d97e12351c066bd"
$ws.Range("B44").Value = "bf0a6fefe132698159578d31a1e6a69c"
$ws.Range("B89").Value = "4c871696d3a94fcb5b8f5daa2f12615b"
$ws.Range("B99").Value = "1b5b59d54e36ae147bfee341efadc733"
$ws.Range("B110").Value = "5bb84315bd01b6d4a1d8ebcbf13f7ce1"
$ws.Range("B121").Value = "69354063445d005bff295dcd858ebfd3"
$ws.Range("B126").Value = "885125c12218fb55a9c17a473ea811f5"
$ws.Range("B133").Value = "c02e3d6b3d2aa91939b0858fb3651db9"
$ws.Range("B136").Value = "b7039b57dbda92005e340256ad999b90"
$ws.Range("B154").Value = "e3bddf8e25ff01c1c35efc6f771131ed"
$ws.Range("B160").Value = "6f8190f43977b1f6fdeacfb05d0efd38"
$ws.Range("B168").Value = "c95d714b63b7a0b2612d21d7a5d379bf"
$ws.Range("B183").Value = "1566ad624c9b683444f8640e7090cefd"
$ws.Range("B191").Value = "97ce9b79e88359ac562527cc9ed8e1a6"
$ws.Range("B198").Value = "9bca5d76692ec5957281453e46621ed6"
$ws.Range("B200").Value = "d5ef55e36803ff9c65c83cdd13fffe52"
$ws.Range("B228").Value = "fe38701a3da4b84079059572acfcc9b3"
$ws.Range("B278").Value = "d6031017e0c1033d48f77d4eae4d20ff"
$ws.Range("B281").Value = "5303e7c7c414586e96e97fca9adc5a1a"
$ws.Range("B302").Value = "c72ab92478c61d71a94c691b800f69f1"
$ws.Range("B335").Value = "57244df0d40fcf86589a51d16474ea7d"
$ws.Range("B464").Value = "88ca15026fa327f90edcf2607339c165"
$ws.Range("B465").Value = "7788fa9a9646e7159463bd9b2733690f"
$ws.Range("B485").Value = "12145fb009908848a5a850e0c8fcda8e"
$ws.Range("B508").Value = "6a55751d6462bd11b65b7440271838e8"
$ws.Range("B542").Value = "0b24743b7eaec31b65d235f0d4706c47"
$ws.Range("B555").Value = "d801b8e81876e7c4a64433dfd4dc2b7b"
$ws.Range("B558").Value = "cb211322d39ea5dcae043e1ec1002c9b"
$ws.Range("B561").Value = "755096d6796763048bca42a3ada26c55"
$ws.Range("B580").Value = "521ce29e8304ca26acab34907e3d08da"
$ws.Range("B592").Value = "5d2a86836ac6c080f9dadf3e71dd678a"
$ws.Range("B637").Value = "2a019fc9d68f80c9cbbe7cd8a452d16d"
$ws.Range("B673").Value = "142e844414e2ca2b6173ed72bee7eb75"
$ws.Range("B708").Value = "12e5dbeb119384264be0298d3ffb04dd"
$ws.Range("B711").Value = "85376c330cb8c179172eb8012e4289fc"
$ws.Range("B712").Value = "c3305368066951b035b3eec49bbfc9ce"
$ws.Range("B723").Value = "3d55dde6eea0e77c61e852a4347905de"
$ws.Range("B750").Value = "bebe597650251d7dc4b5abfc624cebb2"
$ws.Range("B776").Value = "800e6222377b60a6266c2d4131c982a2"
$ws.Range("B819").Value = "3c28e2b2ee8006a3717bce3480372475"
$ws.Range("B823").Value = "a675002953b99d10ce3cc6692c676267"
$ws.Range("B833").Value = "65efb3004a8be5e6b626de0b8267fb17"
$ws.Range("B835").Value = "f10e040faebcd89ecad4e85e77ff55a1"
$ws.Range("B838").Value = "10e0d3fcba82c94ccc94802d6c5c9179"
$ws.Range("B882").Value = "c9c849f03081bb7a17b5eba5feebb7ea"
